$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

# Row 1, columns A:E keep their old "counts_rN" header positions but now hold
# numeric-looking text values (stored as text, same as the original inlineStr cells).
$ws.Range("A1:E1").NumberFormat = "@"
$ws.Range("A1").Value = "0.016650717703349284"
$ws.Range("B1").Value = "0.00616267942583732"
$ws.Range("C1").Value = "0.003827751196172249"
$ws.Range("D1").Value = "0.0022870813397129187"
$ws.Range("E1").Value = "0.0018883572567783093"

# F1/G1 become real numbers (30 and 1) instead of the old duplicated header text.
$ws.Cells.Item(1, 6).Value = 30
$ws.Cells.Item(1, 7).Value = 1

# Everything else that used to live in this sheet (old F1:J1 header duplicates and
# the whole second data row A2:J2) is gone - drop the now-unused columns/row.
$ws.Range("H1:J1").EntireColumn.Delete()
$ws.Range("A2:G2").EntireRow.Delete()
